$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row "5,000" (table row 2): Mantel r -0.024 -> 0.056 ; p 0.295 -> 0.047 (now bold)
$t.Cell(2, 3).Range.Text = "0.056"
$pCell1 = $t.Cell(2, 4)
$pCell1.Range.Text = "0.047"
$pCell1.Range.Font.Bold = $true

# Row "15,000" (table row 3): Mantel r -0.001 -> 0.013 ; p 0.589 -> 0.361
$t.Cell(3, 3).Range.Text = "0.013"
$t.Cell(3, 4).Range.Text = "0.361"

# Row "25,000" (table row 4): Mantel r 0.045 -> 0.032 ; p 0.15 -> 0.208
$t.Cell(4, 3).Range.Text = "0.032"
$t.Cell(4, 4).Range.Text = "0.208"
